$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 81: IdCardNo changes (student accommodation request upload) ---
$ws.Range("C81").Value = 456382

# --- Append new row 82 with the uploaded ID-card / accommodation request data ---
$ws.Range("A82").Value = 35765
$ws.Range("A82").NumberFormat = "m/d/yy"

$ws.Range("B82").Value = "Romania"
$ws.Range("C82").Value = "TC 419786"
$ws.Range("D82").Value = "SPCLEP Tulcea"

$ws.Range("E82").Value = 42539
$ws.Range("E82").NumberFormat = "m/d/yy"

$ws.Range("F82").Value = "Tulcea"
$ws.Range("G82").Value = "Tulcea"
$ws.Range("H82").Value = "Str. Sabinelor"
$ws.Range("I82").Value = "N"

# Move selection to the newly added row, matching the workbook's last edit state
$null = $ws.Range("I82").Select()
